$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Nedas_J (sheet1): fix Week 5 total and append Week 6
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Nedas_J")
$ws1.Range("B5").Value = 0.5625
$ws1.Range("A6").Value = "Week 6"
$ws1.Range("B6").Value = 0.74305555555555547
$ws1.Range("B6").NumberFormat = $ws1.Range("B5").NumberFormat
$ws1.Columns.Item(2).ColumnWidth = 11.17

# ---------------------------------------------------------------------------
# Adomas_J (sheet2): append Week 6
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Adomas_J")
$ws2.Range("A6").Value = "Week 6"
$ws2.Range("B6").Value = 0.46527777777777773
$ws2.Range("B6").NumberFormat = $ws2.Range("B5").NumberFormat
$ws2.Columns.Item(2).ColumnWidth = 11.17

# ---------------------------------------------------------------------------
# Aiste_G (sheet3): append Week 6
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Aistė_G")
$ws3.Range("A6").Value = "Week 6"
$ws3.Range("B6").Value = 0.41319444444444442
$ws3.Range("B6").NumberFormat = $ws3.Range("B5").NumberFormat
$ws3.Columns.Item(2).ColumnWidth = 11.17

# ---------------------------------------------------------------------------
# Gabrielius_D (sheet4): append Week 6
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Gabrielius_D")
$ws4.Range("A6").Value = "Week 6"
$ws4.Range("B6").Value = 0.54166666666666663
$ws4.Range("B6").NumberFormat = $ws4.Range("B5").NumberFormat
$ws4.Columns.Item(2).ColumnWidth = 11.17

# ---------------------------------------------------------------------------
# Overview (sheet5): add Week 6 rows and extend the totals formula
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Overview")

$ws5.Range("A11").Value = "Week 6 (working in a group)"
$ws5.Range("B11").Value = 0
$ws5.Range("B11").NumberFormat = $ws5.Range("B9").NumberFormat

$ws5.Range("A12").Value = "Week 6 (total working hours)"
$ws5.Range("B12").NumberFormat = "[hh]:mm"
$ws5.Range("B12").Font.Color = 4473924
$ws5.Range("B12").Formula = "=SUM(Nedas_J!B6, Adomas_J!B6, Aistė_G!B6, Gabrielius_D!B6) + 4 * B11"

$ws5.Range("E1").Formula = "=SUM(B2+B4+B6+B8+B10+B12)"

# ---------------------------------------------------------------------------
# Selection / view bookkeeping to mirror the final interactive state
# ---------------------------------------------------------------------------
[void]$ws1.Range("B7").Select()
[void]$ws3.Range("B7").Select()
[void]$ws4.Range("P16").Select()
[void]$ws5.Range("B12").Select()

[void]$ws2.Activate()
[void]$ws2.Range("B7").Select()
